# Regen sval data to filter save games
# Update the numeric stat columns (B:G) for rows 2-6 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 3.286832544864788
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 10.19245300693656
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 15.88780690183548

# Row 3
$ws.Range("B3").Value = 0.003208871385164791
$ws.Range("C3").Value = 0.002571899574220771
$ws.Range("D3").Value = 22.3905356188092
$ws.Range("E3").Value = 10.19245300693656
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 32.58876939670514

# Row 4
$ws.Range("B4").Value = 3.286832544864788
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 3.537761648806719
$ws.Range("E4").Value = 0.4942365360607697
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 8.974608811992548

# Row 5
$ws.Range("B5").Value = 3.286832544864788
$ws.Range("C5").Value = 1.655778082260271
$ws.Range("D5").Value = 0.7527432677738641
$ws.Range("E5").Value = 0.4942365360607697
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.189590430959694

# Row 6
$ws.Range("B6").Value = 3.286832544864788
$ws.Range("C6").Value = 1.655778082260271
$ws.Range("D6").Value = 3.537761648806719
$ws.Range("E6").Value = 0.4942365360607697
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 8.974608811992548
